$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 22.63168697499412
$ws.Cells.Item(2, 2).Value = 18.76749301780098
$ws.Cells.Item(2, 3).Value = 26.82010428337896
$ws.Cells.Item(3, 1).Value = 23.9291292432585
$ws.Cells.Item(3, 2).Value = 20.29274859721846
$ws.Cells.Item(3, 3).Value = 27.41812932442699
$ws.Cells.Item(4, 1).Value = 24.05824535099941
$ws.Cells.Item(4, 2).Value = 21.44331589271246
$ws.Cells.Item(4, 3).Value = 27.05438509023186
$ws.Cells.Item(5, 1).Value = 29.07014635957465
$ws.Cells.Item(5, 2).Value = 25.4740593381565
$ws.Cells.Item(5, 3).Value = 33.16134600433896
$ws.Cells.Item(6, 1).Value = 29.87922956369371
$ws.Cells.Item(6, 2).Value = 27.00699985357457
$ws.Cells.Item(6, 3).Value = 32.76927469992776
$ws.Cells.Item(7, 1).Value = 20.64756509430535
$ws.Cells.Item(7, 2).Value = 17.55217860838825
$ws.Cells.Item(7, 3).Value = 24.27586267261828
$ws.Cells.Item(8, 1).Value = 10.75844751802228
$ws.Cells.Item(8, 2).Value = 8.302945801369095
$ws.Cells.Item(8, 3).Value = 13.92987479262075
$ws.Cells.Item(9, 1).Value = 27.83321323252759
$ws.Cells.Item(9, 2).Value = 25.26047765574528
$ws.Cells.Item(9, 3).Value = 30.23759511701108
$ws.Cells.Item(10, 1).Value = 10.87703822765267
$ws.Cells.Item(10, 2).Value = 8.44906911318478
$ws.Cells.Item(10, 3).Value = 14.0337702350185
$ws.Cells.Item(11, 1).Value = 18.67513664405169
$ws.Cells.Item(11, 2).Value = 14.85533363793444
$ws.Cells.Item(11, 3).Value = 22.76638835411477
$ws.Cells.Item(12, 1).Value = 29.38820248727315
$ws.Cells.Item(12, 2).Value = 24.81155639585389
$ws.Cells.Item(12, 3).Value = 33.71311768442493
$ws.Cells.Item(13, 1).Value = 12.85019403546008
$ws.Cells.Item(13, 2).Value = 9.927210507023608
$ws.Cells.Item(13, 3).Value = 16.79280734028303
$ws.Cells.Item(14, 1).Value = 15.20698784370924
$ws.Cells.Item(14, 2).Value = 11.8699624862558
$ws.Cells.Item(14, 3).Value = 19.36086861627459
$ws.Cells.Item(15, 1).Value = 34.02952023242469
$ws.Cells.Item(15, 2).Value = 30.19900118338121
$ws.Cells.Item(15, 3).Value = 37.71982665077127
$ws.Cells.Item(16, 1).Value = 29.91210730467125
$ws.Cells.Item(16, 2).Value = 27.00373809344991
$ws.Cells.Item(16, 3).Value = 32.81325350821284
$ws.Cells.Item(17, 1).Value = 29.25923641245671
$ws.Cells.Item(17, 2).Value = 24.7759776147792
$ws.Cells.Item(17, 3).Value = 33.64522838471533
$ws.Cells.Item(18, 1).Value = 24.08092913765668
$ws.Cells.Item(18, 2).Value = 21.17246752779685
$ws.Cells.Item(18, 3).Value = 27.27915226984385
$ws.Cells.Item(19, 1).Value = 19.23717449044537
$ws.Cells.Item(19, 2).Value = 15.02546725528521
$ws.Cells.Item(19, 3).Value = 24.0324729664443
$ws.Cells.Item(20, 1).Value = 33.58282665725635
$ws.Cells.Item(20, 2).Value = 30.35645023648164
$ws.Cells.Item(20, 3).Value = 37.35221047801411
$ws.Cells.Item(21, 1).Value = 33.63411332327293
$ws.Cells.Item(21, 2).Value = 30.53194336042426
$ws.Cells.Item(21, 3).Value = 36.9885524846296
$ws.Cells.Item(22, 1).Value = 11.17762605625933
$ws.Cells.Item(22, 2).Value = 8.552930182844548
$ws.Cells.Item(22, 3).Value = 14.84245644360468
$ws.Cells.Item(23, 1).Value = 36.41951290086661
$ws.Cells.Item(23, 2).Value = 31.91708899405405
$ws.Cells.Item(23, 3).Value = 40.6865490696918
$ws.Cells.Item(24, 1).Value = 29.91210730467125
$ws.Cells.Item(24, 2).Value = 27.00373809344991
$ws.Cells.Item(24, 3).Value = 32.81325350821284
$ws.Cells.Item(25, 1).Value = 25.6467943267526
$ws.Cells.Item(25, 2).Value = 23.39700306552193
$ws.Cells.Item(25, 3).Value = 28.02226308304904
$ws.Cells.Item(26, 1).Value = 33.72926832786016
$ws.Cells.Item(26, 2).Value = 30.28277324877227
$ws.Cells.Item(26, 3).Value = 37.69655044683951
$ws.Cells.Item(27, 1).Value = 12.99654016223697
$ws.Cells.Item(27, 2).Value = 10.06953269281656
$ws.Cells.Item(27, 3).Value = 16.94029208518673
$ws.Cells.Item(28, 1).Value = 21.87374011067277
$ws.Cells.Item(28, 2).Value = 18.33541881696286
$ws.Cells.Item(28, 3).Value = 25.44041734057405
$ws.Cells.Item(29, 1).Value = 38.68427897489019
$ws.Cells.Item(29, 2).Value = 34.22288902865593
$ws.Cells.Item(29, 3).Value = 43.67495022865955
$ws.Cells.Item(30, 1).Value = 27.97011253211019
$ws.Cells.Item(30, 2).Value = 24.33469417568543
$ws.Cells.Item(30, 3).Value = 31.90126978571163
$ws.Cells.Item(31, 1).Value = 15.11295949299416
$ws.Cells.Item(31, 2).Value = 11.78271528915791
$ws.Cells.Item(31, 3).Value = 18.92372541955596
$ws.Cells.Item(32, 1).Value = 12.52084751499195
$ws.Cells.Item(32, 2).Value = 9.94172212282353
$ws.Cells.Item(32, 3).Value = 15.98907820951386
$ws.Cells.Item(33, 1).Value = 15.19582404933427
$ws.Cells.Item(33, 2).Value = 11.7969406100275
$ws.Cells.Item(33, 3).Value = 19.70321645529663
$ws.Cells.Item(34, 1).Value = 24.41274060533717
$ws.Cells.Item(34, 2).Value = 19.66814385510655
$ws.Cells.Item(34, 3).Value = 29.64879392129155
$ws.Cells.Item(35, 1).Value = 31.75495392168817
$ws.Cells.Item(35, 2).Value = 28.58804703320543
$ws.Cells.Item(35, 3).Value = 35.14843967722261
$ws.Cells.Item(36, 1).Value = 20.68221701682943
$ws.Cells.Item(36, 2).Value = 17.74660429177204
$ws.Cells.Item(36, 3).Value = 24.2623499044183
$ws.Cells.Item(37, 1).Value = 35.84783873336134
$ws.Cells.Item(37, 2).Value = 32.24034679160788
$ws.Cells.Item(37, 3).Value = 39.79707091160779
$ws.Cells.Item(38, 1).Value = 15.38610751880198
$ws.Cells.Item(38, 2).Value = 11.86865253601307
$ws.Cells.Item(38, 3).Value = 19.72867468903289
$ws.Cells.Item(39, 1).Value = 29.95128764850944
$ws.Cells.Item(39, 2).Value = 25.24264140922784
$ws.Cells.Item(39, 3).Value = 34.99753573798459
$ws.Cells.Item(40, 1).Value = 28.9230656480156
$ws.Cells.Item(40, 2).Value = 25.91867359514685
$ws.Cells.Item(40, 3).Value = 31.94856937653019
$ws.Cells.Item(41, 1).Value = 19.42960012366195
$ws.Cells.Item(41, 2).Value = 16.6252568303216
$ws.Cells.Item(41, 3).Value = 22.29588359257033
$ws.Cells.Item(42, 1).Value = 18.88954197034484
$ws.Cells.Item(42, 2).Value = 15.11506198793133
$ws.Cells.Item(42, 3).Value = 23.10173592932824
$ws.Cells.Item(43, 1).Value = 27.96892363483179
$ws.Cells.Item(43, 2).Value = 24.33138010007175
$ws.Cells.Item(43, 3).Value = 31.88639132550494
$ws.Cells.Item(44, 1).Value = 16.71932389291783
$ws.Cells.Item(44, 2).Value = 13.49960903421216
$ws.Cells.Item(44, 3).Value = 19.96915165063418
$ws.Cells.Item(45, 1).Value = 19.2720134975041
$ws.Cells.Item(45, 2).Value = 14.98444537320336
$ws.Cells.Item(45, 3).Value = 24.96814550047718
$ws.Cells.Item(46, 1).Value = 14.64207264562138
$ws.Cells.Item(46, 2).Value = 11.57119326252334
$ws.Cells.Item(46, 3).Value = 18.32088945309552
$ws.Cells.Item(47, 1).Value = 11.27704811928698
$ws.Cells.Item(47, 2).Value = 8.580010904353708
$ws.Cells.Item(47, 3).Value = 14.85744134629589
$ws.Cells.Item(48, 1).Value = 29.49057575149001
$ws.Cells.Item(48, 2).Value = 24.92816727440868
$ws.Cells.Item(48, 3).Value = 33.87207328312751
$ws.Cells.Item(49, 1).Value = 26.44490871440235
$ws.Cells.Item(49, 2).Value = 23.68144243272063
$ws.Cells.Item(49, 3).Value = 29.41225844459385
$ws.Cells.Item(50, 1).Value = 35.7517085984234
$ws.Cells.Item(50, 2).Value = 32.36629911369431
$ws.Cells.Item(50, 3).Value = 39.62250277001835
$ws.Cells.Item(51, 1).Value = 23.91066259617116
$ws.Cells.Item(51, 2).Value = 20.30485073151545
$ws.Cells.Item(51, 3).Value = 27.36937313437732
$ws.Cells.Item(52, 1).Value = 14.21377721366141
$ws.Cells.Item(52, 2).Value = 11.28959671924621
$ws.Cells.Item(52, 3).Value = 17.50566386619605
$ws.Cells.Item(53, 1).Value = 30.02268471943589
$ws.Cells.Item(53, 2).Value = 25.37093613750317
$ws.Cells.Item(53, 3).Value = 35.03450433886296
$ws.Cells.Item(54, 1).Value = 34.14038517828982
$ws.Cells.Item(54, 2).Value = 30.11508624487573
$ws.Cells.Item(54, 3).Value = 37.91970386435871
$ws.Cells.Item(55, 1).Value = 33.70310738658326
$ws.Cells.Item(55, 2).Value = 30.28181671661588
$ws.Cells.Item(55, 3).Value = 37.26954445688628
$ws.Cells.Item(56, 1).Value = 15.05148320982531
$ws.Cells.Item(56, 2).Value = 11.84719468012237
$ws.Cells.Item(56, 3).Value = 19.22733715632063
$ws.Cells.Item(57, 1).Value = 12.52084751499195
$ws.Cells.Item(57, 2).Value = 9.94172212282353
$ws.Cells.Item(57, 3).Value = 15.98907820951386
$ws.Cells.Item(58, 1).Value = 35.91748429223218
$ws.Cells.Item(58, 2).Value = 32.16607510048672
$ws.Cells.Item(58, 3).Value = 39.84666853477985
$ws.Cells.Item(59, 1).Value = 18.31593634397251
$ws.Cells.Item(59, 2).Value = 14.61846495679236
$ws.Cells.Item(59, 3).Value = 21.77064600421878
$ws.Cells.Item(60, 1).Value = 36.34774334342065
$ws.Cells.Item(60, 2).Value = 32.14025727686959
$ws.Cells.Item(60, 3).Value = 40.52303681007199
$ws.Cells.Item(61, 1).Value = 15.37933180991141
$ws.Cells.Item(61, 2).Value = 12.09676549995798
$ws.Cells.Item(61, 3).Value = 19.67404631691175
$ws.Cells.Item(62, 1).Value = 18.67335442880588
$ws.Cells.Item(62, 2).Value = 14.85140486282836
$ws.Cells.Item(62, 3).Value = 22.76638171878558
$ws.Cells.Item(63, 1).Value = 15.35343728764527
$ws.Cells.Item(63, 2).Value = 11.91532797148991
$ws.Cells.Item(63, 3).Value = 20.46591280436887
$ws.Cells.Item(64, 1).Value = 36.35895020780764
$ws.Cells.Item(64, 2).Value = 32.1615501406841
$ws.Cells.Item(64, 3).Value = 40.53847083531604
$ws.Cells.Item(65, 1).Value = 12.29966986075552
$ws.Cells.Item(65, 2).Value = 9.61427016828025
$ws.Cells.Item(65, 3).Value = 15.2264181634325
$ws.Cells.Item(66, 1).Value = 24.13958327112559
$ws.Cells.Item(66, 2).Value = 19.59638404769202
$ws.Cells.Item(66, 3).Value = 28.99999639614118
$ws.Cells.Item(67, 1).Value = 23.28129738231659
$ws.Cells.Item(67, 2).Value = 18.84517638789344
$ws.Cells.Item(67, 3).Value = 28.19161947686301
$ws.Cells.Item(68, 1).Value = 31.13392740002993
$ws.Cells.Item(68, 2).Value = 27.89551958980033
$ws.Cells.Item(68, 3).Value = 34.32718624134862
$ws.Cells.Item(69, 1).Value = 11.20292840287522
$ws.Cells.Item(69, 2).Value = 8.596785885128613
$ws.Cells.Item(69, 3).Value = 14.84590070940418
$ws.Cells.Item(70, 1).Value = 28.91337402920517
$ws.Cells.Item(70, 2).Value = 25.98303500150408
$ws.Cells.Item(70, 3).Value = 31.99772596858741
$ws.Cells.Item(71, 1).Value = 21.88469388615509
$ws.Cells.Item(71, 2).Value = 19.61567080243458
$ws.Cells.Item(71, 3).Value = 24.28714492491569
$ws.Cells.Item(72, 1).Value = 25.23942435813418
$ws.Cells.Item(72, 2).Value = 22.14503976893067
$ws.Cells.Item(72, 3).Value = 28.24886759100151
